$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.145.14'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.515.23'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.74'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.507.83'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.189'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.62'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +11.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.597'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.07'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000275'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '683.27'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.074.74'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.73'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.141.10'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.502.31'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.31%  '

$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.39'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.10'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.909'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.54'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -7.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.47'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.70%  '

$ws.Range("E26").Value = '  -4.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.66'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.42'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.81'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.28'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '578.36'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.68'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -12.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.83'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.105'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.41'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.462.31'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0439'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.336'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0700'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -7.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.133'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.94'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.67%  '

$ws.Range("E51").Value = '  -1.51%  '
